# Refresh the "想去人数" (want-to-go count) figures in column F across the
# 展览 / 演出 / 本地生活 / 全部类型 sheets, matching the regenerated gh-pages
# data snapshot (commit 456a3b4). Only column F values change; everything
# else on each row (price, link, cover, etc.) is left untouched.

$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item("展览")
$wsPerformance = $wb.Worksheets.Item("演出")
$wsLocalLife = $wb.Worksheets.Item("本地生活")
$wsAllTypes = $wb.Worksheets.Item("全部类型")

# 展览
$wsExhibition.Range("F3").Value = 1241
$wsExhibition.Range("F4").Value = 12833
$wsExhibition.Range("F5").Value = 731
$wsExhibition.Range("F10").Value = 1868
$wsExhibition.Range("F11").Value = 40
$wsExhibition.Range("F13").Value = 516
$wsExhibition.Range("F15").Value = 125
$wsExhibition.Range("F18").Value = 297
$wsExhibition.Range("F19").Value = 131
$wsExhibition.Range("F20").Value = 127
$wsExhibition.Range("F23").Value = 249
$wsExhibition.Range("F24").Value = 1294

# 演出
$wsPerformance.Range("F6").Value = 158
$wsPerformance.Range("F17").Value = 12

# 本地生活
$wsLocalLife.Range("F3").Value = 4126

# 全部类型
$wsAllTypes.Range("F6").Value = 1241
$wsAllTypes.Range("F7").Value = 12833
$wsAllTypes.Range("F9").Value = 731
$wsAllTypes.Range("F10").Value = 4126
$wsAllTypes.Range("F15").Value = 1868
$wsAllTypes.Range("F16").Value = 40
$wsAllTypes.Range("F18").Value = 516
$wsAllTypes.Range("F21").Value = 158
$wsAllTypes.Range("F22").Value = 158
$wsAllTypes.Range("F24").Value = 125
$wsAllTypes.Range("F32").Value = 297
$wsAllTypes.Range("F33").Value = 131
$wsAllTypes.Range("F34").Value = 127
$wsAllTypes.Range("F40").Value = 249
$wsAllTypes.Range("F41").Value = 1294
$wsAllTypes.Range("F46").Value = 12
